# Applies the "Updated cryptos list" GitHub Actions commit to Sheet1.
# Every D/E (and a couple of B/C) cell below is stored as text in the
# workbook (inlineStr), so each write is forced to stay text: prefix the
# literal with a leading apostrophe (Excel's quote-prefix, so numeric-
# looking strings like "213.96" are not reinterpreted as numbers), then
# restore the cell style Excel nudges to a quote-prefixed variant back to
# "Normal" so formatting is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $origStyle
}

Set-TextValue "D2" "27.056.07"
Set-TextValue "E2" "  -0.18%  "
Set-TextValue "D3" "1.621.94"
Set-TextValue "E3" "  -0.87%  "
Set-TextValue "E4" "  -0.06%  "
Set-TextValue "D5" "213.96"
Set-TextValue "E5" "  -1.31%  "
Set-TextValue "D6" "0.513"
Set-TextValue "E6" "  -0.70%  "
Set-TextValue "E8" "  +0.39%  "
Set-TextValue "E9" "  -1.36%  "
Set-TextValue "D10" "19.93"
Set-TextValue "E10" "  +0.13%  "
Set-TextValue "D11" "0.0841"
Set-TextValue "E11" "  -0.67%  "
Set-TextValue "D12" "1.848.92"
Set-TextValue "E12" "  -0.87%  "
Set-TextValue "D13" "1.626.90"
Set-TextValue "E13" "  -0.54%  "
Set-TextValue "E14" "  -0.01%  "
Set-TextValue "D15" "0.539"
Set-TextValue "D16" "27.052.22"
Set-TextValue "E16" "  -0.23%  "
Set-TextValue "D17" "64.48"
Set-TextValue "E17" "  -3.24%  "
Set-TextValue "D18" "0.0₃0737"
Set-TextValue "E18" "  -0.35%  "
Set-TextValue "D19" "214.57"
Set-TextValue "E19" "  -0.92%  "
Set-TextValue "E20" "  -0.07%  "
Set-TextValue "D21" "6.83"
Set-TextValue "E21" "  -0.04%  "
Set-TextValue "E22" "  -1.77%  "
Set-TextValue "D23" "2.32"
Set-TextValue "E23" "  -7.57%  "
Set-TextValue "E24" "  -0.85%  "
Set-TextValue "D25" "147.48"
Set-TextValue "E25" "  +0.64%  "
Set-TextValue "E26" "  +0.66%  "
Set-TextValue "E27" "  -0.03%  "
Set-TextValue "E28" "  -3.52%  "
Set-TextValue "D29" "15.51"
Set-TextValue "E29" "  -0.93%  "
Set-TextValue "E30" "  +0.70%  "
Set-TextValue "E31" "  -1.16%  "
Set-TextValue "E32" "  -1.88%  "
Set-TextValue "D33" "0.720"
Set-TextValue "E33" "  +32.83%  "
Set-TextValue "E34" "  -0.14%  "
Set-TextValue "D35" "1.337.68"
Set-TextValue "E35" "  +2.79%  "
Set-TextValue "E36" "  -1.01%  "
Set-TextValue "E37" "  -0.61%  "
Set-TextValue "E38" "  -0.35%  "
Set-TextValue "D39" "0.839"
Set-TextValue "E39" "  -1.62%  "
Set-TextValue "E40" "  -0.09%  "
Set-TextValue "E41" "  -0.31%  "
Set-TextValue "E42" "  -1.55%  "
Set-TextValue "D43" "5.34"
Set-TextValue "E43" "  +0.77%  "
Set-TextValue "D44" "63.92"
Set-TextValue "E44" "  +3.67%  "
Set-TextValue "D45" "1.760.37"
Set-TextValue "E45" "  -0.88%  "
Set-TextValue "D46" "89.88"
Set-TextValue "E46" "  -1.46%  "
Set-TextValue "E47" "  +2.54%  "
Set-TextValue "D48" "0.862"
Set-TextValue "E48" "  +29.23%  "
Set-TextValue "B49" "Cronos"
Set-TextValue "C49" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D49" "0.0512"
Set-TextValue "E49" "  -0.07%  "
Set-TextValue "B50" "Algorand"
Set-TextValue "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.0996"
Set-TextValue "E50" "  +4.16%  "
Set-TextValue "D51" "7.56"
Set-TextValue "E51" "  -1.25%  "
